$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "67.640.45"
$ws.Cells.Item(2,5).Value = "  -1.09%  "

$ws.Cells.Item(3,4).Value = "3.785.43"
$ws.Cells.Item(3,5).Value = "  +1.02%  "

$ws.Cells.Item(4,5).Value = "  +0.05%  "

$ws.Cells.Item(5,4).Value = "'594.98"
$ws.Cells.Item(5,5).Value = "  +0.25%  "

$ws.Cells.Item(6,4).Value = "'166.93"
$ws.Cells.Item(6,5).Value = "  +0.36%  "

$ws.Cells.Item(7,4).Value = "3.771.10"
$ws.Cells.Item(7,5).Value = "  +0.66%  "

$ws.Cells.Item(8,5).Value = "  +0.09%  "

$ws.Cells.Item(9,5).Value = "  +0.34%  "

$ws.Cells.Item(10,5).Value = "  -0.20%  "

$ws.Cells.Item(11,4).Value = "'6.29"
$ws.Cells.Item(11,5).Value = "  -2.33%  "

$ws.Cells.Item(12,5).Value = "  -0.21%  "

$ws.Cells.Item(13,4).Value = "'0.0000252"
$ws.Cells.Item(13,5).Value = "  -2.56%  "

$ws.Cells.Item(14,5).Value = "  -0.59%  "

$ws.Cells.Item(15,4).Value = "4.419.95"
$ws.Cells.Item(15,5).Value = "  +1.07%  "

$ws.Cells.Item(16,4).Value = "3.755.65"
$ws.Cells.Item(16,5).Value = "  +0.17%  "

$ws.Cells.Item(17,4).Value = "'18.52"
$ws.Cells.Item(17,5).Value = "  +3.76%  "

$ws.Cells.Item(18,4).Value = "67.612.42"
$ws.Cells.Item(18,5).Value = "  -1.12%  "

$ws.Cells.Item(19,5).Value = "  +0.36%  "

$ws.Cells.Item(20,5).Value = "  -0.10%  "

$ws.Cells.Item(21,4).Value = "'10.04"
$ws.Cells.Item(21,5).Value = "  -5.85%  "

$ws.Cells.Item(22,4).Value = "'459.32"
$ws.Cells.Item(22,5).Value = "  -1.46%  "

$ws.Cells.Item(23,5).Value = "  -0.14%  "

$ws.Cells.Item(24,5).Value = "  +5.98%  "

$ws.Cells.Item(25,4).Value = "'83.40"
$ws.Cells.Item(25,5).Value = "  -0.68%  "

$ws.Cells.Item(26,4).Value = "'12.00"
$ws.Cells.Item(26,5).Value = "  +0.96%  "

$ws.Cells.Item(27,4).Value = "'2.11"
$ws.Cells.Item(27,5).Value = "  -3.09%  "

$ws.Cells.Item(28,5).Value = "  +0.11%  "

$ws.Cells.Item(29,4).Value = "'9.99"
$ws.Cells.Item(29,5).Value = "  -0.67%  "

$ws.Cells.Item(30,2).Value = "PancakeSwap"
$ws.Cells.Item(30,3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(30,4).Value = "'2.77"
$ws.Cells.Item(30,5).Value = "  +0.24%  "

$ws.Cells.Item(31,2).Value = "ImmutableX"
$ws.Cells.Item(31,3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(31,4).Value = "'2.23"
$ws.Cells.Item(31,5).Value = "  +3.16%  "

$ws.Cells.Item(32,2).Value = "NEARProtocol"
$ws.Cells.Item(32,3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(32,4).Value = "'7.20"
$ws.Cells.Item(32,5).Value = "  -1.39%  "

$ws.Cells.Item(33,2).Value = "EthereumClassic"
$ws.Cells.Item(33,3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(33,4).Value = "'29.62"
$ws.Cells.Item(33,5).Value = "  -0.73%  "

$ws.Cells.Item(34,2).Value = "Binance-PegBSC-USD"
$ws.Cells.Item(34,3).Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Cells.Item(34,4).Value = "'1.00"
$ws.Cells.Item(34,5).Value = "  +0.27%  "

$ws.Cells.Item(35,4).Value = "'9.08"
$ws.Cells.Item(35,5).Value = "  -0.94%  "

$ws.Cells.Item(36,2).Value = "Hedera"
$ws.Cells.Item(36,3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(36,4).Value = "'0.100"
$ws.Cells.Item(36,5).Value = "  -0.61%  "

$ws.Cells.Item(37,2).Value = "dogwifhat"
$ws.Cells.Item(37,3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(37,4).Value = "'3.38"
$ws.Cells.Item(37,5).Value = "  -0.37%  "

$ws.Cells.Item(38,2).Value = "Kaspa"
$ws.Cells.Item(38,3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(38,4).Value = "'0.138"
$ws.Cells.Item(38,5).Value = "  -0.22%  "

$ws.Cells.Item(39,2).Value = "Mantle"
$ws.Cells.Item(39,3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(39,4).Value = "'0.994"
$ws.Cells.Item(39,5).Value = "  -0.40%  "

$ws.Cells.Item(40,2).Value = "Filecoin"
$ws.Cells.Item(40,3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(40,4).Value = "'5.76"
$ws.Cells.Item(40,5).Value = "  -0.32%  "

$ws.Cells.Item(41,2).Value = "FirstDigitalUSD"
$ws.Cells.Item(41,3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(41,4).Value = "'1.00"
$ws.Cells.Item(41,5).Value = "  +0.07%  "

$ws.Cells.Item(42,2).Value = "USDe"
$ws.Cells.Item(42,3).Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Cells.Item(42,5).Value = "  -0.01%  "

$ws.Cells.Item(43,2).Value = "Arweave"
$ws.Cells.Item(43,3).Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Cells.Item(43,4).Value = "'45.28"
$ws.Cells.Item(43,5).Value = "  +3.33%  "

$ws.Cells.Item(44,2).Value = "OKB"
$ws.Cells.Item(44,3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(44,4).Value = "'48.15"
$ws.Cells.Item(44,5).Value = "  +3.20%  "

$ws.Cells.Item(45,2).Value = "TheGraph"
$ws.Cells.Item(45,3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Cells.Item(45,4).Value = "'0.299"
$ws.Cells.Item(45,5).Value = "  -0.95%  "

$ws.Cells.Item(46,2).Value = "Monero"
$ws.Cells.Item(46,3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(46,4).Value = "'149.86"
$ws.Cells.Item(46,5).Value = "  +3.68%  "

$ws.Cells.Item(47,2).Value = "Cosmos"
$ws.Cells.Item(47,3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(47,4).Value = "'8.31"
$ws.Cells.Item(47,5).Value = "  -1.80%  "

$ws.Cells.Item(48,2).Value = "Bittensor"
$ws.Cells.Item(48,3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(48,4).Value = "'394.00"
$ws.Cells.Item(48,5).Value = "  +0.63%  "

$ws.Cells.Item(49,2).Value = "EnergySwap"
$ws.Cells.Item(49,3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(49,4).Value = "'26.65"
$ws.Cells.Item(49,5).Value = "  +5.76%  "

$ws.Cells.Item(50,2).Value = "Stacks"
$ws.Cells.Item(50,3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(50,4).Value = "'1.81"
$ws.Cells.Item(50,5).Value = "  -5.18%  "

$ws.Cells.Item(51,2).Value = "Maker"
$ws.Cells.Item(51,3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(51,4).Value = "2.715.56"
$ws.Cells.Item(51,5).Value = "  -1.19%  "
